$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: new task entry ("3") with a message note in column B.
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = " Testing message"

# Matches the saved selection in the target workbook.
$ws.Range("B5").Select()
